# Apply benchmark update changes to the BENCHMARK sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Row 2
$ws.Range("H2").Value = "15 TL - 15 TL"

# Row 3
$ws.Range("F3").Value = ""

# Row 4
$ws.Range("F4").Value = ""

# Row 5
$ws.Range("F5").Value = ""

# Row 6
$ws.Range("D6").Value = "8.300,01 TL - 199,41 TL"

# Row 7
$ws.Range("H7").Value = "%3,09"

# Row 8
$ws.Range("F8").Value = ""

# Row 9
$ws.Range("F9").Value = ""

# Row 10
$ws.Range("F10").Value = ""

# Row 12
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"

# Row 13
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 11.380 TL"
$ws.Range("F13").Value = ""

# Row 14
$ws.Range("D14").Value = "3.500 TL - 13.500 TL"
$ws.Range("F14").Value = ""
